# Tutorial 6 solution update: change the date separators from "/" to "-"
# for the attendance log dates (rows 3-21), and mark the first two dates
# (28-07-2022 and 01-08-2022) as counted in Total Attendance Count (D)
# and Invalid (G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateText($cellAddress, $text) {
    $c = $ws.Range($cellAddress)
    # Force a text number format first so Excel does not reinterpret the
    # dash-separated string as a real date value, then restore the default
    # "Normal" style so no explicit style index is left on the cell.
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-DateText "A3"  "28-07-2022"
Set-DateText "A4"  "01-08-2022"
Set-DateText "A5"  "04-08-2022"
Set-DateText "A6"  "08-08-2022"
Set-DateText "A7"  "11-08-2022"
Set-DateText "A8"  "15-08-2022"
Set-DateText "A9"  "18-08-2022"
Set-DateText "A10" "22-08-2022"
Set-DateText "A11" "25-08-2022"
Set-DateText "A12" "29-08-2022"
Set-DateText "A13" "01-09-2022"
Set-DateText "A14" "05-09-2022"
Set-DateText "A15" "08-09-2022"
Set-DateText "A16" "12-09-2022"
Set-DateText "A17" "15-09-2022"
Set-DateText "A18" "19-09-2022"
Set-DateText "A19" "22-09-2022"
Set-DateText "A20" "26-09-2022"
Set-DateText "A21" "29-09-2022"

# Rows 3 and 4 also get their Total Attendance Count (D) and Invalid (G)
# counters bumped from 0 to 1
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("D4").Value = 1
$ws.Range("G4").Value = 1
